$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 7:159 down to 8:160 by copying row 6 through row 159 into
# rows 7 through 160. This both pushes the existing data down one row
# and seeds the new row 7 with formatting copied from row 6 (mirroring
# Excel's native "insert row, format from row above" behavior).
$src = $ws.Range("A6:I159")
$dst = $ws.Range("A7:I160")
$src.Copy($dst)

# Fill in the new agenda entry on row 7. Columns that use a "quote
# prefix" text style (e.g. the ID column) need a leading apostrophe so
# the runtime keeps treating them as explicit text instead of dropping
# the quote-prefix formatting.
$ws.Range("A7").Value = "Pedro"
$ws.Range("B7").Value = "'2261"
$ws.Range("C7").Value = "Escola CNA"
$ws.Range("D7").Value = "Atualização no sisetma de alarmes do cliente."
$ws.Range("E7").Value = ""
$ws.Range("G7").Value = "Pendente"
$ws.Range("I7").Value = ""

# Match the saved selection from the edited workbook.
$ws.Range("H7").Select()
